# #5: insurance, claim, debt, investment done
# Fill out the 保險 (Insurance, sheet index 5) and 債務 (Debt, sheet index 6)
# sheets: turn the accidental "row1 == row2" duplicate header into a real
# header row, and append the standard trailing metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that the other sheets already have.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 5: 保險 (Insurance)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

# Build proper column headers in B1:D1 (was a duplicate of row 2's data).
$ws5.Range("B1").Value = "company"
$ws5.Range("C1").Value = "name"
$ws5.Range("D1").Value = "owner"

# New trailing header cells E1:K1 - copy formatting (style s=1) from B1 first.
$ws5.Range("B1").Copy()
$ws5.Range("E1:K1").PasteSpecial(-4122)
$ws5.Range("E1").Value = "property_category"
$ws5.Range("F1").Value = "category"
$ws5.Range("G1").Value = "date"
$ws5.Range("H1").Value = "legislator_name"
$ws5.Range("I1").Value = "legislator_id"
$ws5.Range("J1").Value = "source_file"
$ws5.Range("K1").Value = "index"

# New trailing data cells E2:K5 - copy formatting (style s=2) from B2 first.
$ws5.Range("B2").Copy()
$ws5.Range("E2:K5").PasteSpecial(-4122)

$insRows = @(2, 3, 4, 5)
$insIndex = @(78, 79, 80, 81)
for ($i = 0; $i -lt $insRows.Length; $i++) {
    $r = $insRows[$i]
    $ws5.Range("E$r").Value = "insurance"
    $ws5.Range("F$r").Value = "normal"
    $ws5.Range("G$r").Value = "2013-11-12"
    $ws5.Range("H$r").Value = "王育敏"
    $ws5.Range("I$r").Value = 1728
    $ws5.Range("J$r").Value = "tmped871"
    $ws5.Range("K$r").Value = $insIndex[$i]
}

# ---------------------------------------------------------------------
# Sheet 6: 債務 (Debt)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# Build proper column headers in B1:G1 (was a duplicate of row 2's data).
$ws6.Range("B1").Value = "species"
$ws6.Range("C1").Value = "debtor"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "total"
$ws6.Range("F1").Value = "register_date"
$ws6.Range("G1").Value = "register_reason"

# New trailing header cells H1:N1 - copy formatting (style s=1) from B1 first.
$ws6.Range("B1").Copy()
$ws6.Range("H1:N1").PasteSpecial(-4122)
$ws6.Range("H1").Value = "property_category"
$ws6.Range("I1").Value = "category"
$ws6.Range("J1").Value = "date"
$ws6.Range("K1").Value = "legislator_name"
$ws6.Range("L1").Value = "legislator_id"
$ws6.Range("M1").Value = "source_file"
$ws6.Range("N1").Value = "index"

# New trailing data cells H2:N2 - copy formatting (style s=2) from B2 first.
$ws6.Range("B2").Copy()
$ws6.Range("H2:N2").PasteSpecial(-4122)
$ws6.Range("H2").Value = "debt"
$ws6.Range("I2").Value = "normal"
$ws6.Range("J2").Value = "2013-11-12"
$ws6.Range("K2").Value = "王育敏"
$ws6.Range("L2").Value = 1728
$ws6.Range("M2").Value = "tmped871"
$ws6.Range("N2").Value = 91
